# Apply the "updated board data" changes to the Property Tycoon board data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unhide all the previously-hidden "action" rows on the board.
$hiddenRows = @(5, 7, 9, 10, 12, 15, 17, 20, 22, 25, 27, 30, 33, 35, 38, 40, 41, 43)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $false
}

# 2. Give space #2 (row 6) an unimproved rent value in column I.
$ws.Cells.Item(6, 9).Value = 2

# 3. Rename the "Purple" colour group to "Deep blue" for Skywalker Drive / Tesla Power Co.
$ws.Range("D42").Value = "Deep blue"
$ws.Range("D44").Value = "Deep blue"

# 4. Update the house/hotel-cost notes table to match the renamed colour group.
$ws.Range("H51").Value = "Green, Deep blue"

# 5. Move the active selection to I8, matching the refreshed view.
$ws.Range("I8").Select() | Out-Null
